{"js": "// Apply the LOQ4071.docx content edit: several paragraphs across the\n// document have their text replaced with new literal text (a reshuffle of\n// course-catalog copy), as described by the unified diff.\n//\n// We locate each paragraph by searching for its current (before-edit)\n// text and replace just that text, which preserves paragraph/run\n// formatting (pStyle, bold run-properties, <w:br/> line breaks, etc.).\n\nconst replacements = [\n  {\n    // Objetivos\n    find: \"Complementar a forma\u00e7\u00e3o multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, t\u00f3picos atuais e relevantes sobre gest\u00e3o da qualidade.\",\n    replace: \"A definir, de acordo com o t\u00f3pico programado.\",\n  },\n  {\n    // Docente(s) Respons\u00e1vel(eis)\n    find: \"5840535 - Messias Borges Silva\",\n    replace: \"Complementar a forma\u00e7\u00e3o multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, t\u00f3picos atuais e relevantes sobre gest\u00e3o da qualidade.\",\n  },\n  {\n    // Programa resumido\n    find: \"A definir, de acordo com o t\u00f3pico programado.\",\n    replace: \"O conte\u00fado desta disciplina ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares a forma\u00e7\u00e3o de um profissional de Engenharia.\",\n  },\n  {\n    // Programa\n    find: \"O conte\u00fado desta disciplina ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares a forma\u00e7\u00e3o de um profissional de Engenharia.\",\n    replace: \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\",\n  },\n  {\n    // Avalia\u00e7\u00e3o / M\u00e9todo\n    find: \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\",\n    replace: \"Provas e trabalhos.\",\n  },\n  {\n    // Avalia\u00e7\u00e3o / Crit\u00e9rio\n    find: \"Provas e trabalhos.\",\n    replace: \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\",\n  },\n  {\n    // Avalia\u00e7\u00e3o / Norma de recupera\u00e7\u00e3o \u2014 becomes two lines separated by a\n    // manual line break; \"\\v\" (vertical tab) is how Office.js represents a\n    // <w:br/> soft line break inside inserted text.\n    find: \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\",\n    replace: \"Textos fornecidos pelo professor da disciplina\\vArtigos extra\u00eddos de revistas especializadas na \u00e1rea de gest\u00e3o e produ\u00e7\u00e3o.\",\n  },\n  {\n    // Bibliografia \u2014 two lines collapse into a single line of text.\n    find: \"Textos fornecidos pelo professor da disciplina\\vArtigos extra\u00eddos de revistas especializadas na \u00e1rea de gest\u00e3o e produ\u00e7\u00e3o.\",\n    replace: \"5840535 - Messias Borges Silva\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  // Word.js search supports \"\\v\" (vertical tab) in the search string to\n  // match across a <w:br/> line break, so the full multi-line text (when\n  // present) can be matched and replaced as a single range.\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  // Use the last match: once earlier replacements land new text into the\n  // document, a short search string could coincidentally match more than\n  // once, but the paragraph we still need to edit is always the most\n  // recently-untouched (latest, in document order) occurrence.\n  const target = results.items[results.items.length - 1];\n  target.insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the LOQ4071.docx content edit: several paragraphs across the\n# document have their text replaced with new literal text (a reshuffle of\n# course-catalog copy), as described by the unified diff.\n#\n# Each target paragraph is located by its (stable, 1-based) position in\n# $d.Paragraphs, so replacements never depend on other, not-yet-applied\n# replacements and there is no ambiguity from duplicate text appearing\n# elsewhere in the document.\n\n$d = $word.ActiveDocument\n\n# --- Simple whole-paragraph text swaps -------------------------------\n# These paragraphs contain a single run with no special character\n# formatting, so replacing the paragraph range's Text wholesale is safe\n# and keeps the paragraph mark / style intact.\n\n# Objetivos\n$d.Paragraphs(6).Range.Text = \"A definir, de acordo com o t\u00f3pico programado.\"\n\n# Docente(s) Respons\u00e1vel(eis)\n$d.Paragraphs(8).Range.Text = \"Complementar a forma\u00e7\u00e3o multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, t\u00f3picos atuais e relevantes sobre gest\u00e3o da qualidade.\"\n\n# Programa resumido\n$d.Paragraphs(10).Range.Text = \"O conte\u00fado desta disciplina ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares a forma\u00e7\u00e3o de um profissional de Engenharia.\"\n\n# Programa\n$d.Paragraphs(12).Range.Text = \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\"\n\n# --- Avalia\u00e7\u00e3o paragraph: M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o ----\n# This paragraph mixes bold \"label\" runs with plain \"value\" runs and\n# manual line breaks, so each value is replaced in place with a Find\n# scoped to the paragraph's own Range (never touching the bold labels or\n# the neighbouring values).\n\n$avaliacao = $d.Paragraphs(14).Range\n\n$avaliacao.Find.Text = \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\"\n$avaliacao.Find.Replacement.Text = \"Provas e trabalhos.\"\n$avaliacao.Find.Execute($avaliacao.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $avaliacao.Find.Replacement.Text, 2) | Out-Null\n\n$avaliacao.Find.Text = \"Provas e trabalhos.\"\n$avaliacao.Find.Replacement.Text = \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\"\n$avaliacao.Find.Execute($avaliacao.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $avaliacao.Find.Replacement.Text, 2) | Out-Null\n\n# \"^l\" in the replacement text is Word's manual-line-break code (<w:br/>),\n# splitting this value into two lines within the same run set.\n$avaliacao.Find.Text = \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\"\n$avaliacao.Find.Replacement.Text = \"Textos fornecidos pelo professor da disciplina^lArtigos extra\u00eddos de revistas especializadas na \u00e1rea de gest\u00e3o e produ\u00e7\u00e3o.\"\n$avaliacao.Find.Execute($avaliacao.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $avaliacao.Find.Replacement.Text, 2) | Out-Null\n\n# Bibliografia (two lines separated by a manual line break collapse into\n# a single line of plain text)\n$d.Paragraphs(16).Range.Text = \"5840535 - Messias Borges Silva\"\n"}
